$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the existing row 459. Excel shifts row 459
# (and everything below it, through the former row 543) down by one,
# pushing the former last data row (543) to row 544. This also mirrors the
# D-column date style (s="2") from the surrounding rows onto the new row.
$ws.Rows.Item(459).Insert()

# Populate the newly inserted row 459 with the new weekly price entry.
$ws.Range("A459").Value = 4
$ws.Range("B459").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C459").Value = "Los Lagos"
$ws.Range("D459").Value = 45209
$ws.Range("E459").Value = 10
$ws.Range("F459").Value = 100112045
$ws.Range("G459").Value = "Zapallo"
$ws.Range("H459").Value = "Paine"
$ws.Range("I459").Value = "1a (guarda)"
$ws.Range("J459").Value = 1200
$ws.Range("K459").Value = 800
$ws.Range("L459").Value = 800
$ws.Range("M459").Value = 800
$ws.Range("N459").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O459").Value = "Región de O'Higgins"
$ws.Range("P459").Value = 800
$ws.Range("Q459").Value = 1
$ws.Range("R459").Value = "Hortaliza"
